# Apply updates described by the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": order in row 15 (Order ID 5, Sagar Borse) ---
$orders = $wb.Worksheets.Item("All Orders")
$orders.Range("H15").Value = "CANCELLED"   # Status: NEW -> CANCELLED
$orders.Range("M15").Value = "test order"  # Cancel Reason: "" -> "test order"

# --- Sheet "Daily Summary": totals for 2026-01-13 (row 4) ---
$summary = $wb.Worksheets.Item("Daily Summary")
$summary.Range("D4").Value = 10   # Cancelled: 9 -> 10
$summary.Range("E4").Value = 50   # Revenue: 80 -> 50
$summary.Range("G4").Value = 50   # Pending: 80 -> 50
